# Update "想去人数" (number of people wanting to go) counts on the
# 展览 sheet and the mirrored 全部类型 (all types) aggregate sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F3").Value  = 1350
$wsExhibit.Range("F5").Value  = 1024
$wsExhibit.Range("F6").Value  = 10580
$wsExhibit.Range("F11").Value = 688
$wsExhibit.Range("F13").Value = 12423

# 全部类型 sheet updates (same events, different row offsets)
$wsAll.Range("F4").Value  = 1350
$wsAll.Range("F6").Value  = 1024
$wsAll.Range("F7").Value  = 10580
$wsAll.Range("F12").Value = 688
$wsAll.Range("F14").Value = 12423
